# edit.ps1 — applies the "add 2022-Q1 data" change:
#   * the worksheet that used to be "总计" (the running summary) is renamed to
#     "2022-Q1" and repopulated with that quarter's per-fund holdings detail
#   * a brand-new "总计" worksheet is appended, holding the quarterly summary
#     rows (now including the new 2022-Q1 summary line)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Data
# ---------------------------------------------------------------------------

# Per-fund holdings detail for 2022-Q1.
# Columns: 基金代码, 基金名称, 基金规模, 股票总仓位, 仓位占比, 持有市值(亿元), 仓位排名
$data2022Q1 = @(
    @("005267", "嘉实价值精选股票", "65.04", "92.05", "5.54", "3.6032", 8),
    @("166002", "中欧新蓝筹混合 -A", "129.77", "77.81", "2.68", "3.4778", 7),
    @("011264", "中欧新趋势混合（LOF）X", "89.97", "87.17", "3.17", "2.8520", 6),
    @("166001", "中欧新趋势混合(LOF) -A", "89.97", "87.17", "3.17", "2.8520", 6),
    @("001881", "中欧新趋势混合(LOF) -E", "89.97", "87.17", "3.17", "2.8520", 6),
    @("012533", "嘉实价值驱动一年持有期混合型证券投资基金A", "50.84", "91.01", "4.99", "2.5369", 9),
    @("005233", "广发睿毅领先混合", "40.39", "63.34", "5.75", "2.3224", 3),
    @("001117", "中欧精选灵活配置定期开放混合A", "59.37", "86.95", "3.65", "2.1670", 5),
    @("001890", "中欧精选灵活配置定期开放混合E", "59.37", "86.95", "3.65", "2.1670", 5),
    @("011518", "嘉实价值臻选混合型证券投资基金", "35.83", "89.89", "4.77", "1.7091", 10),
    @("166023", "中欧瑞丰灵活配置混合（LOF）A", "32.40", "85.04", "4.69", "1.5196", 2),
    @("070019", "嘉实价值优势混合", "27.51", "93.14", "4.78", "1.3150", 10),
    @("009909", "嘉实动力先锋混合A", "29.43", "90.99", "4.18", "1.2302", 6),
    @("001763", "广发多策略灵活配置混合", "20.27", "69.36", "5.65", "1.1453", 5),
    @("001718", "工银瑞信物流产业股票", "32.39", "86.93", "3.19", "1.0332", 10),
    @("012647", "中欧洞见一年持有混合", "33.02", "65.46", "2.82", "0.9312", 7),
    @("320001", "诺安平衡混合", "12.68", "73.14", "7.17", "0.9092", 3),
    @("011643", "嘉实时代先锋三年持有期混合型证券投资基金A", "19.08", "91.59", "3.42", "0.6525", 10),
    @("530003", "建信优选成长混合A", "15.75", "75.39", "3.50", "0.5512", 8),
    @("960028", "建信优选成长混合H", "15.75", "75.39", "3.50", "0.5512", 8),
    @("001044", "嘉实新消费股票", "8.92", "80.25", "5.97", "0.5325", 5),
    @("004355", "嘉实丰和灵活配置混合", "9.22", "85.59", "4.28", "0.3946", 9),
    @("004237", "中欧新蓝筹混合 -C", "9.82", "77.81", "2.68", "0.2632", 7),
    @("320018", "诺安新动力混合", "3.36", "79.24", "5.00", "0.1680", 4),
    @("012534", "嘉实价值驱动一年持有期混合型证券投资基金C", "3.20", "91.01", "4.99", "0.1597", 9),
    @("005335", "浙商全景消费混合", "2.30", "93.36", "5.75", "0.1322", 10),
    @("009910", "嘉实动力先锋混合C", "2.95", "90.99", "4.18", "0.1233", 6),
    @("004740", "中欧瑞丰灵活配置混合（LOF）C", "1.28", "85.04", "4.69", "0.0600", 2),
    @("168101", "九泰锐智事件驱动混合（LOF）", "0.73", "89.22", "7.88", "0.0575", 7),
    @("530012", "建信积极配置混合", "1.52", "51.74", "2.79", "0.0424", 7),
    @("011644", "嘉实时代先锋三年持有期混合型证券投资基金C", "1.12", "91.59", "3.42", "0.0383", 10),
    @("001885", "中欧新蓝筹混合 -E", "1.41", "77.81", "2.68", "0.0378", 7),
    @("008135", "华宸未来价值先锋混合", "0.20", "86.99", "7.88", "0.0158", 3),
    @("009700", "长江添利混合A", "1.90", "20.46", "0.75", "0.0142", 10),
    @("009701", "长江添利混合C", "1.40", "20.46", "0.75", "0.0105", 10),
    @("005901", "诺安汇利灵活配置混合A", "0.08", "86.88", "6.59", "0.0053", 5),
    @("005902", "诺安汇利灵活配置混合C", "0.02", "86.88", "6.59", "0.0013", 5)
)

# Quarterly summary rows (newest first).
# Columns: 日期, 持有数量(只), 持有市值(亿元)
$dataTotal = @(
    @("2022-Q1", 37, 38.43),
    @("2021-Q4", 56, 48.32),
    @("2021-Q3", 46, 51.56),
    @("2021-Q2", 37, 33.67),
    @("2021-Q1", 62, 34.88),
    @("2020-Q4", 46, 18.17)
)

# ---------------------------------------------------------------------------
# Step 1: turn the old "总计" sheet into the new "2022-Q1" detail sheet
# ---------------------------------------------------------------------------

$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"
$q1.Cells.Clear()

# Borrow the header / index-column cell formatting used by the other
# per-quarter sheets so the new sheet matches their look.
$styleSrc = $wb.Worksheets.Item("2021-Q4")

$q1.Cells.Item(1, 2).Value = "基金代码"
$q1.Cells.Item(1, 3).Value = "基金名称"
$q1.Cells.Item(1, 4).Value = "基金规模"
$q1.Cells.Item(1, 5).Value = "股票总仓位"
$q1.Cells.Item(1, 6).Value = "仓位占比"
$q1.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q1.Cells.Item(1, 8).Value = "仓位排名"

$styleSrc.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

$q1RowCount = $data2022Q1.Count
$q1LastRow = $q1RowCount + 1

# Fund code / name / size / position columns are stored as text (codes have
# leading zeros, sizes keep a fixed 2-4 decimal format) — force text so the
# values round-trip exactly instead of being parsed as numbers.
$q1.Range("B2:G$q1LastRow").NumberFormat = "@"

for ($i = 0; $i -lt $data2022Q1.Count; $i++) {
    $r = $i + 2
    $row = $data2022Q1[$i]
    $q1.Cells.Item($r, 1).Value = $i
    $q1.Cells.Item($r, 2).Value = $row[0]
    $q1.Cells.Item($r, 3).Value = $row[1]
    $q1.Cells.Item($r, 4).Value = $row[2]
    $q1.Cells.Item($r, 5).Value = $row[3]
    $q1.Cells.Item($r, 6).Value = $row[4]
    $q1.Cells.Item($r, 7).Value = $row[5]
    $q1.Cells.Item($r, 8).Value = $row[6]
}

$styleSrc.Range("A2:A$q1LastRow").Copy()
$q1.Range("A2:A$q1LastRow").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Step 2: append a fresh "总计" sheet after "2022-Q1" with the updated
#          quarter-over-quarter summary (incl. the new 2022-Q1 row)
# ---------------------------------------------------------------------------

$total = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q1)
$total.Name = "总计"

$total.Cells.Item(1, 2).Value = "日期"
$total.Cells.Item(1, 3).Value = "持有数量(只)"
$total.Cells.Item(1, 4).Value = "持有市值(亿元)"

$styleSrc.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)

for ($i = 0; $i -lt $dataTotal.Count; $i++) {
    $r = $i + 2
    $row = $dataTotal[$i]
    $total.Cells.Item($r, 1).Value = $i
    $total.Cells.Item($r, 2).Value = $row[0]
    $total.Cells.Item($r, 3).Value = $row[1]
    $total.Cells.Item($r, 4).Value = $row[2]
}

$totalLastRow = $dataTotal.Count + 1
$styleSrc.Range("A2:A$totalLastRow").Copy()
$total.Range("A2:A$totalLastRow").PasteSpecial(-4122)

Write-Output "2022-Q1 sheet inserted; 总计 sheet rebuilt"
